# Update column C (Fitness) values on Sheet1 according to the new run data.
# Rows 2-252 correspond to Generation 0-250 (column B), Run 12 (column A, unchanged).
# Column C values change in blocks as follows:
#   rows 2-7   (Generation 0-5)    -> 8540
#   rows 8-11  (Generation 6-9)    -> 8501
#   rows 12-16 (Generation 10-14)  -> 8499
#   rows 17-24 (Generation 15-22)  -> 7318
#   rows 25-252 (Generation 23-250) -> 7310

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-ColumnCBlock($startRow, $endRow, $value) {
    for ($r = $startRow; $r -le $endRow; $r++) {
        $ws.Cells.Item($r, 3).Value = $value
    }
}

Set-ColumnCBlock 2 7 8540
Set-ColumnCBlock 8 11 8501
Set-ColumnCBlock 12 16 8499
Set-ColumnCBlock 17 24 7318
Set-ColumnCBlock 25 252 7310
